$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values as per the diff
$ws.Range("E3").Value = 12.802
$ws.Range("A9").Value = -20.912
$ws.Range("A18").Value = -21.825
$ws.Range("A20").Value = -21.738
